# Edit: on the F1-score slide, change the denominator of the F1 formula
# from "Precision*Recall" to "Precision+Recall" (the "/ Precision*Recall"
# part becomes "/ Precision+Recall"), leaving the numerator
# "Precision*Recall" untouched.
#
# i.e. "F1 = Precision*Recall / Precision*Recall " (sic, as authored)
#   -> "F1 = Precision*Recall / Precision+Recall "

$p = $ppt.ActivePresentation

# Locate the shape holding the F1 formula text, searching every slide so the
# script keeps working even if slide/shape ordering differs from what we
# expect.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t.Contains("Precision*Recall / Precision*Recall")) {
                $targetShape = $shp
            }
        }
    }
}

if ($targetShape -eq $null) {
    Write-Output "Target shape (F1 formula) not found!"
} else {
    $tr = $targetShape.TextFrame.TextRange

    # Replace " / Precision*" (the run right before the second "Recall")
    # with " / " - keeps the existing run boundaries/formatting for the
    # surrounding text intact.
    $full = $tr.Text
    $searchStr = " / Precision*"
    $idx = $full.IndexOf($searchStr)

    if ($idx -ge 0) {
        $run = $tr.Characters($idx + 1, $searchStr.Length)
        $run.Text = " / "

        # The trailing "Recall" run (originally the 2nd "Recall") now
        # becomes "Precision+Recall".
        $newFull = $tr.Text
        $recallIdx = $newFull.IndexOf("Recall", $idx)
        $recallRun = $tr.Characters($recallIdx + 1, "Recall".Length)
        $recallRun.Text = "Precision+Recall"
    }

    Write-Output "Updated F1 formula text: [$($tr.Text)]"
}
